$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 37 with data filled in from ESR data (mirrors MoH daily report columns)
$ws.Range("A37").Value = 43924
$ws.Range("A37").NumberFormat = $ws.Range("A36").NumberFormat

$ws.Range("B37").Value = 49
$ws.Range("C37").Value = 772
$ws.Range("D37").Value = 22
$ws.Range("E37").Value = 96
$ws.Range("F37").Value = 71
$ws.Range("G37").Value = 868
$ws.Range("H37").Value = 11
$ws.Range("I37").Value = 103
$ws.Range("J37").Value = 13
$ws.Range("L37").Value = 1
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 1
$ws.Range("O37").Value = 425
$ws.Range("P37").Value = 286
$ws.Range("Q37").Value = 148
$ws.Range("R37").Value = 9
$ws.Range("S37").Value = 868
